{"js": "// The document contains a single results table (\"environmental contrasts\")\n// with a header row followed by section-header rows (1 cell) and data rows\n// (4 cells: District | Estimate (95% CI) | p value ... actually: label |\n// district | estimate | p value). Every data row's last (\"p value\") cell is\n// being re-reported with one extra significant digit plus a trailing \"*\"\n// (and values that used to print as \"0\"/\"<2e-16\" become \"<1e-04*\").\n//\n// We target each p-value cell by its absolute (row, column) position in the\n// table, which is the most precise way to reproduce the diff exactly since\n// several old values (\"0\", \"0.021\", \"0.000\", \"<2e-16\") are not unique and\n// map to different replacements depending on which row they are in.\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// rowIndex -> new p-value text (0-based row index into the table, matching\n// Word.Table.getCell()).\nconst pValueUpdates = {\n  2: \"<1e-04*\",\n  4: \"<1e-04*\",\n  5: \"0.0267*\",\n  6: \"0.0210*\",\n  7: \"<1e-04*\",\n  9: \"<1e-04*\",\n  10: \"0.0376*\",\n  11: \"0.0209*\",\n  13: \"<1e-04*\",\n  14: \"0.0084*\",\n  15: \"0.0029*\",\n  17: \"<1e-04*\",\n  19: \"<1e-04*\",\n  20: \"0.0171*\",\n  22: \"<1e-04*\",\n  23: \"0.0003*\",\n  24: \"0.0003*\",\n  26: \"<1e-04*\",\n};\n\nconst pValueColumn = 3; // 0-based: label | district | estimate (95% CI) | p value\n\nfor (const [rowIndexStr, newText] of Object.entries(pValueUpdates)) {\n  const rowIndex = Number(rowIndexStr);\n  const cell = table.getCell(rowIndex, pValueColumn);\n  cell.value = newText;\n}\nawait context.sync();\n\n// styles.xml changes ------------------------------------------------------\n\n// 1) Delete the unused \"Abstract Title\" paragraph style entirely.\nconst abstractTitleStyle = context.document.getStyles().getByNameOrNullObject(\"Abstract Title\");\nabstractTitleStyle.load(\"nameLocal\");\nawait context.sync();\nif (!abstractTitleStyle.isNullObject) {\n  abstractTitleStyle.delete();\n}\n\n// 2) Delete the unused \"Footnote Block Text\" paragraph style entirely.\nconst footnoteBlockTextStyle = context.document.getStyles().getByNameOrNullObject(\"Footnote Block Text\");\nfootnoteBlockTextStyle.load(\"nameLocal\");\nawait context.sync();\nif (!footnoteBlockTextStyle.isNullObject) {\n  footnoteBlockTextStyle.delete();\n}\n\n// 3) \"Abstract\" style: space-before goes from 100 (twips => 5pt) to 300\n//    (twips => 15pt) so it matches space-after.\nconst abstractStyle = context.document.getStyles().getByNameOrNullObject(\"Abstract\");\nabstractStyle.load(\"nameLocal\");\nawait context.sync();\nif (!abstractStyle.isNullObject) {\n  abstractStyle.paragraphFormat.spaceBefore = 15; // points (OOXML w:before is in twentieths of a point)\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single results table (\"environmental contrasts\")\n# with a header row followed by section-header rows (1 cell) and data rows\n# (4 cells: label | district | estimate (95% CI) | p value). Every data\n# row's last (\"p value\") cell is being re-reported with one extra\n# significant digit plus a trailing \"*\" (and values that used to print as\n# \"0\"/\"<2e-16\" become \"<1e-04*\").\n#\n# We target each p-value cell by its absolute (row, column) position in the\n# table (Word COM is 1-indexed), which is the most precise way to reproduce\n# the diff exactly since several old values (\"0\", \"0.021\", \"0.000\",\n# \"<2e-16\") are not unique and map to different replacements depending on\n# which row they are in.\n\n$doc = $word.ActiveDocument\n$table = $doc.Tables.Item(1)\n\n# 1-based row index (Word COM) -> new p-value text.\n$pValueUpdates = [ordered]@{\n    3  = \"<1e-04*\"\n    5  = \"<1e-04*\"\n    6  = \"0.0267*\"\n    7  = \"0.0210*\"\n    8  = \"<1e-04*\"\n    10 = \"<1e-04*\"\n    11 = \"0.0376*\"\n    12 = \"0.0209*\"\n    14 = \"<1e-04*\"\n    15 = \"0.0084*\"\n    16 = \"0.0029*\"\n    18 = \"<1e-04*\"\n    20 = \"<1e-04*\"\n    21 = \"0.0171*\"\n    23 = \"<1e-04*\"\n    24 = \"0.0003*\"\n    25 = \"0.0003*\"\n    27 = \"<1e-04*\"\n}\n\n$pValueColumn = 4  # 1-based: label | district | estimate (95% CI) | p value\n\nforeach ($rowIndex in $pValueUpdates.Keys) {\n    $cell = $table.Cell($rowIndex, $pValueColumn)\n    $cell.Range.Text = $pValueUpdates[$rowIndex]\n}\n\n# styles.xml changes ------------------------------------------------------\n# NB: look styles up individually via Styles.Item(name) (wrapped in\n# try/catch) rather than `foreach ($s in $doc.Styles)` \u2014 iterating the live\n# collection after a deletion re-indexes it and blows up.\n\n# 1) Delete the unused \"Abstract Title\" paragraph style entirely.\ntry {\n    $abstractTitleStyle = $doc.Styles.Item(\"Abstract Title\")\n    $abstractTitleStyle.Delete()\n} catch {\n}\n\n# 2) Delete the unused \"Footnote Block Text\" paragraph style entirely.\ntry {\n    $footnoteBlockTextStyle = $doc.Styles.Item(\"Footnote Block Text\")\n    $footnoteBlockTextStyle.Delete()\n} catch {\n}\n\n# 3) \"Abstract\" style: space-before goes from 100 (twips => 5pt) to 300\n#    (twips => 15pt) so it matches space-after.\n$abstractStyle = $doc.Styles.Item(\"Abstract\")\n$abstractStyle.ParagraphFormat.SpaceBefore = 15\n"}
